$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 21 (2025Q3) metrics as per latest data refresh
$ws.Range("C21").Value = 246
$ws.Range("D21").Value = 218
$ws.Range("E21").Value = 28
$ws.Range("F21").Value = 62.46418338108882
